$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number (45202 -> 2023-10-04 serial 45203)
# for every data row (rows 2 through 135). Increment it by one day for all of them.
$ws.Range("C2:C135").Value = 45203
